$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Beta)
$ws.Range("F2").Value = 40.39741067362148
$ws.Range("G2").Value = 40.02245244293899
$ws.Range("H2").Value = 40.77365862722088
$ws.Range("I2").Value = 0.000765824070251822
$ws.Range("J2").Value = 0.0007163308392964616
$ws.Range("K2").Value = 0.0008583964211809417
$ws.Range("L2").Value = 0.05788682845466459
$ws.Range("M2").Value = 0.05752239324992253
$ws.Range("N2").Value = 0.05825498948212531

# Row 3 (Gamma)
$ws.Range("F3").Value = 0.00001392089462884584
$ws.Range("G3").Value = 0.0000000004514459885090415
$ws.Range("H3").Value = 0.00003969226573169534
$ws.Range("I3").Value = 0.00001214812706397338
$ws.Range("J3").Value = 0.0000000004141827379701685
$ws.Range("K3").Value = 0.00003451564012716813
$ws.Range("L3").Value = 0.00001433227727611874
$ws.Range("M3").Value = 0.0000000004779920793195923
$ws.Range("N3").Value = 0.00004084575416167441

# Row 4 (Beta + Gamma)
$ws.Range("F4").Value = 40.39742459451611
$ws.Range("G4").Value = 40.02245244339043
$ws.Range("H4").Value = 40.77369831948661
$ws.Range("I4").Value = 0.0007779721973157954
$ws.Range("J4").Value = 0.0007163312534791996
$ws.Range("K4").Value = 0.0008929120613081097
$ws.Range("L4").Value = 0.0579011607319407
$ws.Range("M4").Value = 0.05752239372791462
$ws.Range("N4").Value = 0.05829583523628697
